# Update the logged Timestamp (B2) on the "Student Details" sheet to the
# latest recorded attempt, matching the newly appended shared-string entries:
#   03/14/2020 16:43:19
#   03/14/2020 21:48:06
#   03/14/2020 21:48:25
#   03/14/2020 21:48:58   <- newest entry, now referenced by B2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student Details")

$ws.Range("B2").Value = "03/14/2020 21:48:58"
